$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.083.89'
$ws.Range('E2').Value = '  +3.70%  '

$ws.Range('D3').Value = '3.197.00'
$ws.Range('E3').Value = '  +2.03%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '537.68'
$ws.Range('E5').Value = '  +0.49%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.16'
$ws.Range('E6').Value = '  +4.51%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'

$ws.Range('E8').Value = '  +2.27%  '

$ws.Range('E9').Value = '  -0.24%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.113'
$ws.Range('E10').Value = '  +4.52%  '

$ws.Range('E11').Value = '  +2.46%  '

$ws.Range('D12').Value = '3.745.17'
$ws.Range('E12').Value = '  +2.13%  '

$ws.Range('E13').Value = '  -0.87%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.12'
$ws.Range('E14').Value = '  +1.58%  '

$ws.Range('E15').Value = '  +3.69%  '

$ws.Range('D16').Value = '60.096.13'
$ws.Range('E16').Value = '  +3.58%  '

$ws.Range('D17').Value = '3.202.07'
$ws.Range('E17').Value = '  +2.57%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.20'
$ws.Range('E18').Value = '  -0.53%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.13'
$ws.Range('E19').Value = '  +1.54%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.35'
$ws.Range('E20').Value = '  +2.57%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.93'
$ws.Range('E21').Value = '  +1.91%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.24%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.529'
$ws.Range('E23').Value = '  +3.42%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.04'
$ws.Range('E24').Value = '  +0.22%  '

$ws.Range('E25').Value = '  +3.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.82'
$ws.Range('E26').Value = '  +14.52%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').Value = '  -0.28%  '

$ws.Range('D28').Value = '0.0₃0902'
$ws.Range('E28').Value = '  +1.57%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.91'
$ws.Range('E29').Value = '  +1.85%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.47'
$ws.Range('E30').Value = '  +3.48%  '

$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.18'
$ws.Range('E31').Value = '  -0.17%  '

$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.42'
$ws.Range('E32').Value = '  +4.94%  '

$ws.Range('E33').Value = '  +3.31%  '

$ws.Range('E34').Value = '  +4.69%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '156.77'
$ws.Range('E35').Value = '  -2.48%  '

$ws.Range('E36').Value = '  +1.92%  '

$ws.Range('D37').Value = '2.766.26'
$ws.Range('E37').Value = '  +8.56%  '

$ws.Range('E38').Value = '  +0.66%  '

$ws.Range('E39').Value = '  +5.31%  '

$ws.Range('E40').Value = '  +0.73%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.27'
$ws.Range('E41').Value = '  +1.49%  '

$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.78'
$ws.Range('E42').Value = '  +2.77%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.730'
$ws.Range('E43').Value = '  +4.38%  '

$ws.Range('E44').Value = '  +5.85%  '

$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.01'
$ws.Range('E45').Value = '  +3.28%  '

$ws.Range('B46').Value = 'RenzoRestakedETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D46').Value = '3.240.24'
$ws.Range('E46').Value = '  +2.23%  '

$ws.Range('E47').Value = '  -0.13%  '

$ws.Range('E48').Value = '  +2.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.62'
$ws.Range('E49').Value = '  +2.62%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.793'
$ws.Range('E50').Value = '  +5.73%  '
